$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55, shifting the existing rows 55-84 down to 56-85.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly price record.
$ws.Cells.Item(55, 1).Value  = 11
$ws.Cells.Item(55, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value  = "Bíobío"
$ws.Cells.Item(55, 4).Value  = 44673
$ws.Cells.Item(55, 5).Value  = 8
$ws.Cells.Item(55, 6).Value  = 100112001
$ws.Cells.Item(55, 7).Value  = "Berenjena"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 220
$ws.Cells.Item(55, 11).Value = 6000
$ws.Cells.Item(55, 12).Value = 6500
$ws.Cells.Item(55, 13).Value = 6227
$ws.Cells.Item(55, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(55, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value = 104
$ws.Cells.Item(55, 17).Value = 60
$ws.Cells.Item(55, 18).Value = "Hortaliza"
